$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: refreshed crypto price/volume data,
# plus a few re-ordered rows (XRP/USDC, Maker/RenderToken, Hedera/OKB).

$ws.Range("D2").Value = "67.087.55"
$ws.Range("E2").Value = "  +4.31%  "
$ws.Range("D3").Value = "3.266.32"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'580.41"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").Value = "'177.41"
$ws.Range("E6").Value = "  +3.77%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "3.264.48"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").Value = "'6.76"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("E12").Value = "  +4.61%  "
$ws.Range("D13").Value = "3.834.26"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'28.22"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "67.121.90"
$ws.Range("E16").Value = "  +4.38%  "
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("D18").Value = "3.268.30"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D19").Value = "'5.85"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "'13.46"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "'372.94"
$ws.Range("E21").Value = "  +5.70%  "
$ws.Range("D22").Value = "'7.63"
$ws.Range("E22").Value = "  +6.09%  "
$ws.Range("D23").Value = "'71.65"
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").Value = "3.411.53"
$ws.Range("E26").Value = "  +3.14%  "
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").Value = "'9.86"
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("D32").Value = "'5.64"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'22.71"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("D37").Value = "'166.82"
$ws.Range("E37").Value = "  +7.94%  "
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("D39").Value = "'0.858"
$ws.Range("E39").Value = "  +5.85%  "
$ws.Range("E40").Value = "  +10.92%  "
$ws.Range("D41").Value = "'27.31"
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'6.48"
$ws.Range("E43").Value = "  +6.78%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.745.43"
$ws.Range("E44").Value = "  +4.32%  "
$ws.Range("D45").Value = "'4.37"
$ws.Range("E45").Value = "  +4.58%  "
$ws.Range("D46").Value = "'348.07"
$ws.Range("E46").Value = "  +4.83%  "
$ws.Range("D47").Value = "'25.12"
$ws.Range("E47").Value = "  +5.16%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'40.53"
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0678"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").Value = "'0.0282"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").Value = "'0.103"
$ws.Range("E51").Value = "  +1.96%  "
